# CHECKER DONE week13 end
# Adds two more weekly review columns (N, O) to the tracking sheet,
# mirroring the existing weekly date header + reviewer assignment rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: week-start dates (mirror formatting of the previous date cell, M1) ---
$ws.Range("N1").NumberFormat = "m/d/yy"
$ws.Range("N1").Value() = 42338
$ws.Range("O1").NumberFormat = "m/d/yy"
$ws.Range("O1").Value() = 42345

# --- Row 2: "Рецензент 1" assignments ---
$ws.Range("N2").Value() = "Бурамбекова"
$ws.Range("O2").Value() = "Акимутин"

# --- Row 3: "Рецензент 2" assignments ---
$ws.Range("N3").Value() = "Заварзин"
$ws.Range("O3").Value() = "Асеев"

# --- Column widths for the two newly used columns ---
$ws.Columns("N:N").ColumnWidth = 12.417
$ws.Columns("O:O").ColumnWidth = 9.25

# --- Update selection to follow the last filled header cell, as Excel would ---
$ws.Range("O4").Select() | Out-Null
